$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "coded" (numeric-category) duplicate rows for the country-you-live-in,
# state-you-live-in, country-you-work-in, state-you-work-in, and work-position
# questions are redundant now that those questions are left as free-text ("string")
# -- only the "string" row survives for each. Delete bottom-up so the row numbers
# of rows still to be processed do not shift underneath us.
$ws.Rows(112).Delete()
$ws.Rows(110).Delete()
$ws.Rows(108).Delete()
$ws.Rows(106).Delete()
$ws.Rows(104).Delete()

# "What is your gender?" (string row, now row 101): record the value categorization
# lists used to bucket the free-text answers into male / female / other / nan.
#
# NOTE: each of these three lists is typed with a leading "'" in the source file.
# Excel's cell-input parser always treats *one* leading apostrophe as the literal
# text-prefix escape character (stripped from the stored value, cell flagged
# quotePrefix) no matter how the value is written into the cell. The female-list
# value's stored text itself begins with a literal apostrophe (it is quoting
# `'female'` as the first item), so it is written here with *two* leading quotes
# so that after the escape-stripping exactly one literal apostrophe remains --
# matching the stored text exactly.
$ws.Range("D101").Value = '''Male'', ''male'', ''Male '', ''M'', ''m'', ''man'', ''Cis male'', ''Male.'', ''Male (cis)'', ''Sex is male'', ''cis male'', ''Dude'', "I''m a man why didn''t you make this a drop down question. You should of asked sex? And I would of answered yes please. Seriously how much text can this take? ", ''male '', ''Cis Male'', ''cisdude'', ''cis man'', ''MALE'', ''Man'''
$ws.Range("E101").Value = '''''female'', ''Female'', ''I identify as female.'', ''female '', ''Female assigned at birth '', ''F'', ''Woman'', ''f'', ''Cis female '', ''Female '', ''woman'', ''female/woman'', ''Female (props for making this a freeform field, though)'', '' Female'', ''Cis-woman'', ''Cisgender Female'''
$ws.Range("F101").Value = '''Bigender'', ''non-binary'', ''fm'', ''Transitioned, M2F'', ''Genderfluid (born female)'', ''Other/Transfeminine'', ''Female or Multi-Gender Femme'', ''Androgynous'', ''male 9:1 female, roughly'',  ''Other'', ''nb masculine'', ''none of your business'', ''genderqueer'', ''Human'', ''Genderfluid'',''Enby'', ''Malr'', ''genderqueer woman'', ''mtf'', ''Queer'', ''Agender'', ''Fluid'', ''mail'', ''M|'', ''Male/genderqueer'', ''fem'', ''Nonbinary'', ''human'', ''Unicorn'', ''Male (trans, FtM)'', ''Genderqueer'', ''Genderflux demi-girl'', ''female-bodied; no feelings about gender'', ''AFAB'', ''Transgender woman'' '
$ws.Range("G101").Value = "nan"

# "What is your gender?" (coded row, now row 102): note + category codes.
$ws.Range("C102").Value = 'all nan values are category 3'
$ws.Range("D102").Value = 1
$ws.Range("E102").Value = 2
$ws.Range("F102").Value = 3
$ws.Range("G102").Value = 3

# Country/state questions (now rows 103-106): just needed a note that the raw
# string values are usable as-is.
$ws.Range("C103").Value = 'String Value is fine here'
$ws.Range("C104").Value = 'String Value is fine here'
$ws.Range("C105").Value = 'String Value is fine here'
$ws.Range("C106").Value = 'String Value is fine here'

# Work position question (moved up to row 107 after the row deletions): note that
# it allows multiple answers.
$ws.Range("C107").Value = "multiple answers, could possibly also be categorized"

# "Do you work remotely?" (now rows 108-109): value categorization + codes.
$ws.Range("D108").Value = 'Always'
$ws.Range("E108").Value = "Sometimes"
$ws.Range("F108").Value = "Never"
$ws.Range("D109").Value = 1
$ws.Range("E109").Value = 2
$ws.Range("F109").Value = 3

# Scroll / selection, matching where the author left off editing.
$ws.Application.Goto($ws.Range("A98"), $false)
$ws.Range("F114:F115").Select()

